$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.935.08'
$ws.Range("E2").Value = '  +6.24%  '
$ws.Range("D3").Value = '2.517.63'
$ws.Range("E3").Value = '  +3.64%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.43'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.18'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.525'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +2.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.05'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +7.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0816'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.44'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.19'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("D15").Value = '2.899.26'
$ws.Range("E15").Value = '  +3.33%  '
$ws.Range("D16").Value = '2.498.06'
$ws.Range("E16").Value = '  +2.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.83%  '
$ws.Range("D18").Value = '47.770.41'
$ws.Range("E18").Value = '  +6.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.71'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.60'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.56%  '
$ws.Range("D21").Value = '0.0₃0940'
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.89'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.02'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +5.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.57'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.38'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.87%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.13%  '
$ws.Range("E29").Value = '  +6.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.26'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +7.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.135'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +9.41%  '
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.16'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.39'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0785'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.27%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.96'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.67'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +5.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.00'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.84%  '
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '122.03'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.24'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.27'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0298'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.27%  '
$ws.Range("D45").Value = '1.973.29'
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.00'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.33%  '
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.82'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.23'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.39'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +14.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.50'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.54%  '
